$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new blank row at position 7 (shifts FUCICORT row 7->8, total row 8->9,
#    footer row 9->10).
$ws.Rows("7:7").Insert()

# 2. Copy formatting from the (now-shifted) FUCICORT row 8 onto the new blank row 7 so the
#    new item row matches the existing item-row look (fonts/number formats/borders/fills).
$ws.Range("A8:Q8").Copy()
$ws.Range("A7:Q7").PasteSpecial(-4122)
$ws.Rows("7:7").RowHeight = 25.5

# 3. Re-create the merged cells for the new row 7 (merging doesn't come along with
#    PasteSpecial formats).
$ws.Range("A7:B7").Merge()
$ws.Range("C7:G7").Merge()
$ws.Range("H7:K7").Merge()
$ws.Range("L7:M7").Merge()
$ws.Range("N7:O7").Merge()

# 4. Fill in the new item's data (CERELAC rice, no-milk baby food).
$ws.Range("A7").Value = 2
$ws.Range("C7").Value = "CERELAC رز بدون لبن"
$ws.Range("H7").Value = "0:0"
$ws.Range("L7").Value = "0"
$ws.Range("N7").Value = "40.00"
$ws.Range("P7").Value = "40.0000"
$ws.Range("Q7").Value = "1:0"

# Original first item becomes item #2 now that CERELAC is item #1... actually keep the
# original numbering: CERELAC is listed first (#1) and FUCICORT moves to #2.
$ws.Range("A7").Value = 1
$ws.Range("A8").Value = 2

# 5. Update the running total (now on row 9): 40.00 + 70.00 = 110.
$ws.Range("P9").Value = 110

# 6. Update the footer timestamp (now on row 10).
$ws.Range("A10").Value = "Saturday, 6 September, 2025 9:47 AM"
